$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.325.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.884.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.69%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.88%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.504"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.883.31"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.146"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.432"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.76%  "

$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.74%  "

$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.355.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.370.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.08%  "

$ws.Range("E18").Value = "  -1.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.875.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.43%  "

$ws.Range("E21").Value = "  -2.04%  "

$ws.Range("E22").Value = "  -1.17%  "

$ws.Range("E23").Value = "  -2.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "

$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.03%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("E28").Value = "  -4.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000105"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.99%  "

$ws.Range("E32").Value = "  -7.78%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("E34").Value = "  -1.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.84%  "

$ws.Range("E36").Value = "  -3.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.47%  "

$ws.Range("E40").Value = "  -3.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.267"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.703.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0335"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "341.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.02%  "

$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("E50").Value = "  -1.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.37%  "
